$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to end with a merged "TODO" note row (A23:C23) reminding
# the author to add a test for runtime errors. That reminder is replaced
# by a concrete new backlog row ("var usage tree" / var-name-length checks)
# appended right after the existing "?" / "RT" rows.

# Remove the old merged note row entirely (also drops its merge + style).
$ws.Rows("23").Delete()

# Add the new backlog row 22, matching the layout of rows 20-21
# (A = "?", B = "RT"-style category replaced with "Me", C = description).
$ws.Cells.Item(22, 1).Value = "?"
$ws.Cells.Item(22, 2).Value = "Me"
$ws.Cells.Item(22, 3).Value = "Length of var name checks"

# Column A on this sheet carries an explicit style (quotePrefix) applied
# to every other row in the table; copy it onto the new row's A cell too.
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)  # xlPasteFormats

# Restore the no-selection (non-highlighted) cursor position used after
# the edit in the authored workbook.
$ws.Range("C26").Select()
